$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "37.322.85"
$c.Style = "Normal"
$ws.Range("E2").Value = "  +2.26%  "

# Row 3
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "2.060.47"
$c.Style = "Normal"
$ws.Range("E3").Value = "  +3.46%  "

# Row 4
$ws.Range("E4").Value = "  +0.12%  "

# Row 5
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "234.18"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -0.81%  "

# Row 6
$ws.Range("E6").Value = "  +2.56%  "

# Row 7
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "58.04"
$c.Style = "Normal"
$ws.Range("E7").Value = "  +5.58%  "

# Row 8
$ws.Range("E8").Value = "  +0.05%  "

# Row 9
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.381"
$c.Style = "Normal"
$ws.Range("E9").Value = "  +2.85%  "

# Row 10
$ws.Range("E10").Value = "  +1.92%  "

# Row 11
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.0760"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +1.84%  "

# Row 12
$ws.Range("E12").Value = "  +2.70%  "

# Row 13
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "2.365.89"
$c.Style = "Normal"
$ws.Range("E13").Value = "  +3.71%  "

# Row 14
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "14.58"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +3.02%  "

# Row 15
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "21.13"
$c.Style = "Normal"
$ws.Range("E15").Value = "  +3.22%  "

# Row 16
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "0.773"
$c.Style = "Normal"
$ws.Range("E16").Value = "  +2.34%  "

# Row 17
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "5.17"
$c.Style = "Normal"
$ws.Range("E17").Value = "  +1.96%  "

# Row 18
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "2.094.64"
$c.Style = "Normal"
$ws.Range("E18").Value = "  +5.34%  "

# Row 19
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "37.545.05"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +3.02%  "

# Row 20
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "6.16"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +16.65%  "

# Row 21
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "68.95"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +1.71%  "

# Row 22
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "0.0₃0811"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +0.80%  "

# Row 23
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "226.50"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +2.24%  "

# Row 24
$ws.Range("E24").Value = "  +0.13%  "

# Row 25
$ws.Range("E25").Value = "  +1.35%  "

# Row 26
$ws.Range("E26").Value = "  +1.22%  "

# Row 27
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "165.14"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +1.64%  "

# Row 28
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "1.49"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +11.89%  "

# Row 29
$ws.Range("E29").Value = "  +2.07%  "

# Row 30
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "19.13"
$c.Style = "Normal"
$ws.Range("E30").Value = "  +1.25%  "

# Row 31
$ws.Range("E31").Value = "  -1.63%  "

# Row 32
$ws.Range("E32").Value = "  +1.55%  "

# Row 33
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "4.49"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +2.83%  "

# Row 34
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "0.0620"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +2.41%  "

# Row 35
$ws.Range("E35").Value = "  +8.98%  "

# Row 36
$ws.Range("E36").Value = "  +6.97%  "

# Row 37
$ws.Range("B37").Value = "WEMIXToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "1.80"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +2.20%  "

# Row 38
$ws.Range("E38").Value = "  +0.09%  "

# Row 39
$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "3.37"
$c.Style = "Normal"
$ws.Range("E39").Value = "  -0.13%  "

# Row 40
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "5.85"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +4.77%  "

# Row 41
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "0.0981"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +4.23%  "

# Row 42
$ws.Range("E42").Value = "  -1.64%  "

# Row 43
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "4.38"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +23.67%  "

# Row 44
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "1.455.47"
$c.Style = "Normal"
$ws.Range("E44").Value = "  +0.18%  "

# Row 45
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "95.46"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +7.04%  "

# Row 46
$ws.Range("E46").Value = "  +3.72%  "

# Row 47
$ws.Range("E47").Value = "  +4.40%  "

# Row 48
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "15.78"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +4.07%  "

# Row 49
$ws.Range("E49").Value = "  +3.33%  "

# Row 50
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "7.24"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +5.47%  "

# Row 51
$ws.Range("E51").Value = "  +2.16%  "

Write-Host "Updated cryptos list"
